# Commit message: "a lot of mess" — tidying up view state:
#  - rename the second sheet ("Sheet1" -> "ממוין מודל ומקור")
#  - change its zoom level to 85% and move the selection to L8
#  - (workbook window-size bookkeeping is host chrome and isn't
#    something this script can influence)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "ממוין מודל ומקור"

$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("L8").Select()
